# Update column G ("K" - strikeouts) with regenerated values computed from
# the source (replacing the old "Strike#" derived figures).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    3  = 3
    4  = 9
    5  = 7
    6  = 7
    7  = 12
    8  = 5
    9  = 9
    10 = 9
    11 = 9
    12 = 11
    13 = 9
    14 = 13
    15 = 12
    16 = 16
    17 = 14
    18 = 9
    19 = 6
    20 = 1
    21 = 5
    22 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
